$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.965.00"
$ws.Range("E2").Value = "  +6.31%  "
$ws.Range("D3").Value = "'1.883.22"
$ws.Range("E3").Value = "  +5.70%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'248.59"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.4976"
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("D8").Value = "'45.90"
$ws.Range("E8").Value = "  +9.45%  "
$ws.Range("D9").Value = "'0.2852"
$ws.Range("E9").Value = "  +7.01%  "
$ws.Range("D10").Value = "'0.06537"
$ws.Range("E10").Value = "  +4.71%  "
$ws.Range("D11").Value = "'1.881.28"
$ws.Range("E11").Value = "  +5.59%  "
$ws.Range("D12").Value = "'17.08"
$ws.Range("E12").Value = "  +3.75%  "
$ws.Range("D13").Value = "'0.07222"
$ws.Range("E13").Value = "  +2.97%  "
$ws.Range("D14").Value = "'0.6630"
$ws.Range("E14").Value = "  +6.04%  "
$ws.Range("D15").Value = "'85.07"
$ws.Range("E15").Value = "  +6.56%  "
$ws.Range("D16").Value = "'4.790"
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("D17").Value = "'29.966.47"
$ws.Range("E17").Value = "  +6.42%  "
$ws.Range("D18").Value = "'0.9992"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'12.87"
$ws.Range("E19").Value = "  +7.02%  "
$ws.Range("D20").Value = "'0.000007499"
$ws.Range("E20").Value = "  +3.96%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "'2.121.30"
$ws.Range("E22").Value = "  +5.61%  "
$ws.Range("D23").Value = "'4.746"
$ws.Range("E23").Value = "  +4.23%  "
$ws.Range("D24").Value = "'5.538"
$ws.Range("E24").Value = "  +6.03%  "
$ws.Range("D25").Value = "'9.001"
$ws.Range("E25").Value = "  +3.07%  "
$ws.Range("D26").Value = "'145.13"
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("D27").Value = "'134.97"
$ws.Range("E27").Value = "  +23.74%  "
$ws.Range("D29").Value = "'1.955"
$ws.Range("E29").Value = "  +5.20%  "
$ws.Range("D30").Value = "'1.376"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").Value = "'4.174"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "'0.08604"
$ws.Range("E32").Value = "  +4.33%  "
$ws.Range("D33").Value = "'3.868"
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("D34").Value = "'0.05107"
$ws.Range("D35").Value = "'1.128"
$ws.Range("E35").Value = "  +5.62%  "
$ws.Range("D36").Value = "'0.6858"
$ws.Range("E36").Value = "  +5.42%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'1.000"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "'2.704"
$ws.Range("E38").Value = "  +3.54%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'2.307"
$ws.Range("E39").Value = "  +13.07%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.747"
$ws.Range("E40").Value = "  +6.20%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.9582"
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.01630"
$ws.Range("E42").Value = "  +5.40%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'6.069"
$ws.Range("E43").Value = "  +2.27%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'104.26"
$ws.Range("E44").Value = "  +4.59%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "'0.4212"
$ws.Range("E46").Value = "  +5.94%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.435"
$ws.Range("E47").Value = "  +3.83%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1251"
$ws.Range("E48").Value = "  +4.24%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05632"
$ws.Range("E49").Value = "  +3.68%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'32.35"
$ws.Range("E50").Value = "  +5.80%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'8.278"
$ws.Range("E51").Value = "  +3.25%  "
